$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.310.43"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "2.950.68"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.39"
$ws.Range("E5").Value = "  -2.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.78"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "2.947.32"
$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  -4.44%  "

$ws.Range("E11").Value = "  -4.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  -3.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.49"
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("E15").Value = "  -1.47%  "

$ws.Range("D16").Value = "65.291.25"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").Value = "3.412.82"
$ws.Range("E17").Value = "  -2.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.01"
$ws.Range("E18").Value = "  +1.07%  "

$ws.Range("D19").Value = "2.949.25"
$ws.Range("E19").Value = "  -2.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.83"
$ws.Range("E20").Value = "  +13.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.29"
$ws.Range("E21").Value = "  -2.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  +1.79%  "

$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.38"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  -2.23%  "

$ws.Range("E26").Value = "  -1.35%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  -5.98%  "

$ws.Range("E29").Value = "  +5.40%  "

$ws.Range("E30").Value = "  -2.42%  "

$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("E32").Value = "  -1.98%  "

$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.12"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.72"
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.16"
$ws.Range("E38").Value = "  -1.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.74"
$ws.Range("E39").Value = "  +2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("E40").Value = "  -7.69%  "

$ws.Range("E41").Value = "  -1.06%  "

$ws.Range("E42").Value = "  -0.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.83"
$ws.Range("E43").Value = "  -5.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.55"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "381.23"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("E46").Value = "  -1.52%  "

$ws.Range("D47").Value = "2.680.82"
$ws.Range("E47").Value = "  -3.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.43"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.60"
$ws.Range("E50").Value = "  -0.89%  "

$ws.Range("E51").Value = "  +1.23%  "
